$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows 98-102 (revised figures from MV -datos-) ---
$ws.Range("N98").Value = 48979
$ws.Range("P98").Value = 48979
$ws.Range("T98").Value = 97911
$ws.Range("V98").Value = 87473

$ws.Range("N99").Value = 47956
$ws.Range("P99").Value = 47956
$ws.Range("T99").Value = 96525
$ws.Range("V99").Value = 86008

$ws.Range("H100").Value = 12484
$ws.Range("J100").Value = 12483
$ws.Range("K100").Value = 742
$ws.Range("M100").Value = 706
$ws.Range("N100").Value = 47915
$ws.Range("P100").Value = 47915
$ws.Range("T100").Value = 94308
$ws.Range("V100").Value = 85490

$ws.Range("H101").Value = 12596
$ws.Range("J101").Value = 12596
$ws.Range("K101").Value = 751
$ws.Range("M101").Value = 714
$ws.Range("N101").Value = 49055
$ws.Range("P101").Value = 49055
$ws.Range("T101").Value = 97441
$ws.Range("V101").Value = 88733

$ws.Range("N102").Value = 48964
$ws.Range("P102").Value = 48962
$ws.Range("T102").Value = 98989
$ws.Range("V102").Value = 90270

# --- New monthly row for 01-06-2021 ---
# Force text storage so the date-like label isn't auto-converted to a date
# serial, then restore the cell to the default (unstyled) look so no stray
# number format lingers on the new cell, matching the rest of column A.
$ws.Range("A103").NumberFormat = "@"
$ws.Range("A103").Value = "01-06-2021"
$ws.Range("A103").Style = "Normal"

$ws.Range("B103").Value = 36834
$ws.Range("C103").Value = 8807
$ws.Range("D103").Value = 28027
$ws.Range("E103").Value = 12
$ws.Range("F103").Value = 12
$ws.Range("G103").Value = 0
$ws.Range("H103").Value = 12328
$ws.Range("I103").Value = 1
$ws.Range("J103").Value = 12328
$ws.Range("K103").Value = 683
$ws.Range("L103").Value = 36
$ws.Range("M103").Value = 648
$ws.Range("N103").Value = 49620
$ws.Range("O103").Value = 2
$ws.Range("P103").Value = 49618
$ws.Range("Q103").Value = 209
$ws.Range("R103").Value = 0
$ws.Range("S103").Value = 209
$ws.Range("T103").Value = 99687
$ws.Range("U103").Value = 8857
$ws.Range("V103").Value = 90830
